# Corrects stock-report rows where the Item Code / Sale Price / Qty / Value
# columns (B, D, E, F, G) had been mis-assigned across neighbouring line
# items for the same product. Each assignment below restores the value
# that belongs on that row per the canonical report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B142").Value = 48654
$ws.Range("E142").Value = 38.26
$ws.Range("F142").Value = -1
$ws.Range("G142").Value = -32.02
$ws.Range("B143").Value = 63902
$ws.Range("E143").Value = 34.04
$ws.Range("F143").Value = 2
$ws.Range("G143").Value = 64.04000000000001
$ws.Range("B154").Value = 53925
$ws.Range("E154").Value = 79.37
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 66.44
$ws.Range("B155").Value = 64350
$ws.Range("E155").Value = 70.63
$ws.Range("F155").Value = 101
$ws.Range("G155").Value = 6710.44
$ws.Range("B156").Value = 57756
$ws.Range("F156").Value = -100
$ws.Range("G156").Value = -6644
$ws.Range("B256").Value = 48719
$ws.Range("E256").Value = 353.35
$ws.Range("F256").Value = -81
$ws.Range("G256").Value = -23955.75
$ws.Range("B257").Value = 64979
$ws.Range("E257").Value = 314.41
$ws.Range("F257").Value = 82
$ws.Range("G257").Value = 24251.5
$ws.Range("B271").Value = 48706
$ws.Range("E271").Value = 39.8
$ws.Range("F271").Value = -144
$ws.Range("G271").Value = -4795.2
$ws.Range("B272").Value = 64973
$ws.Range("E272").Value = 35.4
$ws.Range("F272").Value = 150
$ws.Range("G272").Value = 4995
$ws.Range("B308").Value = 61610
$ws.Range("D308").Value = 102.71
$ws.Range("E308").Value = 122.71
$ws.Range("F308").Value = -58
$ws.Range("G308").Value = -5957.18
$ws.Range("B309").Value = 63565
$ws.Range("E309").Value = 109.19
$ws.Range("F309").Value = 60
$ws.Range("G309").Value = 6162.6
$ws.Range("B310").Value = 57077
$ws.Range("D310").Value = 93.08
$ws.Range("E310").Value = 111.2
$ws.Range("F310").Value = 1
$ws.Range("G310").Value = 93.08
$ws.Range("B338").Value = 63520
$ws.Range("E338").Value = 153.4
$ws.Range("F338").Value = 97
$ws.Range("G338").Value = 13995.16
$ws.Range("B339").Value = 55373
$ws.Range("E339").Value = 163.62
$ws.Range("F339").Value = -94
$ws.Range("G339").Value = -13562.32
$ws.Range("B342").Value = 63571
$ws.Range("E342").Value = 152.53
$ws.Range("F342").Value = 29
$ws.Range("G342").Value = 4160.92
$ws.Range("B343").Value = 63531
$ws.Range("F343").Value = 80
$ws.Range("G343").Value = 11478.4
$ws.Range("B344").Value = 57802
$ws.Range("E344").Value = 162.71
$ws.Range("F344").Value = -79
$ws.Range("G344").Value = -11334.92
$ws.Range("B367").Value = 61605
$ws.Range("E367").Value = 133.78
$ws.Range("F367").Value = -13
$ws.Range("G367").Value = -1455.48
$ws.Range("B368").Value = 63563
$ws.Range("E368").Value = 119.04
$ws.Range("F368").Value = 15
$ws.Range("G368").Value = 1679.4
$ws.Range("B374").Value = 60325
$ws.Range("E374").Value = 151.57
$ws.Range("F374").Value = -102
$ws.Range("G374").Value = -12939.72
$ws.Range("B375").Value = 63560
$ws.Range("E375").Value = 134.87
$ws.Range("F375").Value = 104
$ws.Range("G375").Value = 13193.44
$ws.Range("B381").Value = 62865
$ws.Range("F381").Value = 151
$ws.Range("G381").Value = 12051.31
$ws.Range("B382").Value = 57817
$ws.Range("F382").Value = 3
$ws.Range("G382").Value = 239.43
$ws.Range("B411").Value = 63007
$ws.Range("F411").Value = 984
$ws.Range("G411").Value = 168588.72
$ws.Range("B412").Value = 57856
$ws.Range("F412").Value = 2
$ws.Range("G412").Value = 342.66
$ws.Range("B578").Value = 45695
$ws.Range("E578").Value = 23.58
$ws.Range("F578").Value = -36
$ws.Range("G578").Value = -710.28
$ws.Range("B579").Value = 64915
$ws.Range("E579").Value = 20.98
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 789.2
$ws.Range("B582").Value = 45706
$ws.Range("E582").Value = 23.58
$ws.Range("F582").Value = -202
$ws.Range("G582").Value = -3985.46
$ws.Range("B583").Value = 64922
$ws.Range("E583").Value = 20.98
$ws.Range("F583").Value = 207
$ws.Range("G583").Value = 4084.11
$ws.Range("B585").Value = 64927
$ws.Range("E585").Value = 17.26
$ws.Range("F585").Value = 295
$ws.Range("G585").Value = 4784.9
$ws.Range("B586").Value = 45718
$ws.Range("E586").Value = 19.38
$ws.Range("F586").Value = -294
$ws.Range("G586").Value = -4768.68
$ws.Range("B679").Value = 53319
$ws.Range("E679").Value = 310.64
$ws.Range("F679").Value = -6
$ws.Range("G679").Value = -1643.52
$ws.Range("B680").Value = 64810
$ws.Range("E680").Value = 291.22
$ws.Range("F680").Value = 7
$ws.Range("G680").Value = 1917.44
$ws.Range("B701").Value = 60025
$ws.Range("E701").Value = 37.22
$ws.Range("F701").Value = -98
$ws.Range("G701").Value = -3217.34
$ws.Range("B702").Value = 64833
$ws.Range("E702").Value = 34.9
$ws.Range("F702").Value = 99
$ws.Range("G702").Value = 3250.17
$ws.Range("B712").Value = 60022
$ws.Range("E712").Value = 37.22
$ws.Range("F712").Value = -113
$ws.Range("G712").Value = -3709.79
$ws.Range("B713").Value = 64830
$ws.Range("E713").Value = 34.9
$ws.Range("F713").Value = 117
$ws.Range("G713").Value = 3841.11
$ws.Range("B864").Value = 54751
$ws.Range("E864").Value = 46.34
$ws.Range("F864").Value = -19
$ws.Range("G864").Value = -776.53
$ws.Range("B865").Value = 65079
$ws.Range("F865").Value = 21
$ws.Range("G865").Value = 858.27
